$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header changes ---
# A1 gets new "GC" header
$ws.Range("A1").Value = "GC"
# D1 label renamed (was "replace_pos_<50%")
$ws.Range("D1").Value = "delete/insert"

# New "Simple" block label (set early so shared-string insertion order matches)
$ws.Range("A7").Value = "Simple"

# E1 / F1 labels renamed
$ws.Range("E1").Value = "pos<50%"
$ws.Range("F1").Value = "pos>=50%"

# --- Row 2 (Nodes) values: D/E/F swap+update ---
$ws.Range("D2").Value = 3952
$ws.Range("E2").Value = 522
$ws.Range("F2").Value = 485

# --- Row 3 (CC) values ---
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0.48
$ws.Range("F3").Value = 0.51

# --- Row 4 (ASPL) values ---
$ws.Range("D4").Value = 13.23
$ws.Range("E4").Value = 5.57
$ws.Range("F4").Value = 6.67

# --- Row 5 (S) values ---
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 8.9
$ws.Range("F5").Value = 36.159999999999997

# --- New block starting row 7 ---
$ws.Range("A8").Value = "Nodes"
$ws.Range("B8").Value = 9075
$ws.Range("C8").Value = 7742
$ws.Range("D8").Value = 6248
$ws.Range("E8").Value = 7233
$ws.Range("F8").Value = 6729

$ws.Range("A9").Value = "CC"
$ws.Range("B9").Value = 0.31
$ws.Range("C9").Value = 0.39
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0.53
$ws.Range("F9").Value = 0.63

$ws.Range("A10").Value = "ASPL"
$ws.Range("B10").Value = 6.05
$ws.Range("C10").Value = 5.28
$ws.Range("D10").Value = 13.22
$ws.Range("E10").Value = 5.0199999999999996
$ws.Range("F10").Value = 5.47

# --- Selection / view state update ---
$ws.Range("F22").Select()
